$wb = $excel.ActiveWorkbook

# ===========================================================================
# Sheet "Overview" -- rows re-ordered (ffffd20f672e, ffffff58abeb4c, b7b1e851)
# and the b7b1e851 record's status flips to "In Translation" with a new
# "Latest HO Xliff Generate Date".
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$ws.Range("B2").Value = "e2e\ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$ws.Range("G2").Value = "2017-02-21 04:48:42"

$ws.Range("A3").Value = "ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"
$ws.Range("B3").Value = "e2e\ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"

$ws.Range("A4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$ws.Range("B4").Value = "e2e\b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2017-02-21 04:56:14"

$addrOvB2 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$addrOvB3 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5410a325da766eb8e98c7a977911f7c65626d205/e2e/ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$addrOvB4 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $addrOvB2, [Type]::Missing, [Type]::Missing, "e2e\ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $addrOvB3, [Type]::Missing, [Type]::Missing, "e2e\ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md")
$ws.Hyperlinks.Add($ws.Range("B4"), $addrOvB4, [Type]::Missing, [Type]::Missing, "e2e\b7b1e851-a32e-47b1-9c3f-841bf00595e0.md")

# ===========================================================================
# Sheet "zh-cn" -- same re-ordering + handoff/handback bookkeeping refresh for
# the b7b1e851 row (now last), plus a stale-handback error message.
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$ws.Range("G2").Value = "81fee5a2-f1b9-4318-b938-70094a70c564.ab19f837a08b0166ba47420fc757d1710d134b77.zh-cn.xlf"
$ws.Range("H2").Value = "2017-02-21 04:48:27"
$ws.Range("J2").Value = "81fee5a2-f1b9-4318-b938-70094a70c564.md"
$ws.Range("K2").Value = "81fee5a2-f1b9-4318-b938-70094a70c564.ab19f837a08b0166ba47420fc757d1710d134b77.zh-cn.xlf"
$ws.Range("L2").Value = "2017-02-21 04:49:22"

$ws.Range("A3").Value = "ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"
$ws.Range("F3").Value = "True"
$ws.Range("H3").Value = "2017-02-21 04:48:42"

$ws.Range("A4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.3939f3c8d1bf75def6a2cbe224b2397ec7064d13.zh-cn.xlf"
$ws.Range("H4").Value = "2017-02-21 04:55:57"
$ws.Range("J4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$ws.Range("K4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.3939f3c8d1bf75def6a2cbe224b2397ec7064d13.zh-cn.xlf"
$ws.Range("L4").Value = "2017-02-21 04:54:39"
$ws.Range("R4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/3c6435669bbe22d731a3c108ffa0befa75d1bb30/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md."

$ws.Columns.Item(18).ColumnWidth = 39.14

$addrZhA2 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$addrZhJ2 = "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/34e6f9805f0a28346195a27ad7f2c55e9d517c53/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$addrZhA3 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5410a325da766eb8e98c7a977911f7c65626d205/e2e/ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$addrZhJ3 = "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/9a29a408c1e1e4fab7570928b3eafdb922fb9ca3/e2e/81fee5a2-f1b9-4318-b938-70094a70c564.md"
$addrZhA4 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"
$addrZhJ4 = "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/9a29a408c1e1e4fab7570928b3eafdb922fb9ca3/e2e/81fee5a2-f1b9-4318-b938-70094a70c564.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $addrZhA2, [Type]::Missing, [Type]::Missing, "ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md")
$ws.Hyperlinks.Add($ws.Range("J2"), $addrZhJ2, [Type]::Missing, [Type]::Missing, "81fee5a2-f1b9-4318-b938-70094a70c564.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $addrZhA3, [Type]::Missing, [Type]::Missing, "ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md")
$ws.Hyperlinks.Add($ws.Range("J3"), $addrZhJ3, [Type]::Missing, [Type]::Missing, "81fee5a2-f1b9-4318-b938-70094a70c564.md")
$ws.Hyperlinks.Add($ws.Range("A4"), $addrZhA4, [Type]::Missing, [Type]::Missing, "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md")
$ws.Hyperlinks.Add($ws.Range("J4"), $addrZhJ4, [Type]::Missing, [Type]::Missing, "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md")

# ===========================================================================
# Sheet "de-de" -- mirrors the zh-cn changes, with the de-de specific xlf
# file names / dates.
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$ws.Range("G2").Value = "81fee5a2-f1b9-4318-b938-70094a70c564.ab19f837a08b0166ba47420fc757d1710d134b77.de-de.xlf"
$ws.Range("H2").Value = "2017-02-21 04:48:42"
$ws.Range("J2").Value = "81fee5a2-f1b9-4318-b938-70094a70c564.md"
$ws.Range("K2").Value = "81fee5a2-f1b9-4318-b938-70094a70c564.ab19f837a08b0166ba47420fc757d1710d134b77.de-de.xlf"
$ws.Range("L2").Value = "2017-02-21 04:49:45"

$ws.Range("A3").Value = "ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"
$ws.Range("F3").Value = "True"

$ws.Range("A4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.3939f3c8d1bf75def6a2cbe224b2397ec7064d13.de-de.xlf"
$ws.Range("H4").Value = "2017-02-21 04:56:14"
$ws.Range("J4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$ws.Range("K4").Value = "b7b1e851-a32e-47b1-9c3f-841bf00595e0.3939f3c8d1bf75def6a2cbe224b2397ec7064d13.de-de.xlf"
$ws.Range("L4").Value = "2017-02-21 04:55:02"
$ws.Range("R4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/3c6435669bbe22d731a3c108ffa0befa75d1bb30/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md."

$ws.Columns.Item(18).ColumnWidth = 39.14

$addrDeA2 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$addrDeJ2 = "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/8aece7ab7165b4c0daf4ffe6c0d406e74fb3c8a6/e2e/b7b1e851-a32e-47b1-9c3f-841bf00595e0.md"
$addrDeA3 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5410a325da766eb8e98c7a977911f7c65626d205/e2e/ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md"
$addrDeJ3 = "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/c3bde7029989d05f826ede1f4e74cd8c494d70db/e2e/81fee5a2-f1b9-4318-b938-70094a70c564.md"
$addrDeA4 = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/db4dce94dc46fa9e2d3d6ceec324daae4d194347/e2e/ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md"
$addrDeJ4 = "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/c3bde7029989d05f826ede1f4e74cd8c494d70db/e2e/81fee5a2-f1b9-4318-b938-70094a70c564.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $addrDeA2, [Type]::Missing, [Type]::Missing, "ffffd20f672e-5e6c-49cc-98ee-295e83f53e11.md")
$ws.Hyperlinks.Add($ws.Range("J2"), $addrDeJ2, [Type]::Missing, [Type]::Missing, "81fee5a2-f1b9-4318-b938-70094a70c564.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $addrDeA3, [Type]::Missing, [Type]::Missing, "ffffff58abeb4c-4e19-45c6-b592-3d8ef01f04dc.md")
$ws.Hyperlinks.Add($ws.Range("J3"), $addrDeJ3, [Type]::Missing, [Type]::Missing, "81fee5a2-f1b9-4318-b938-70094a70c564.md")
$ws.Hyperlinks.Add($ws.Range("A4"), $addrDeA4, [Type]::Missing, [Type]::Missing, "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md")
$ws.Hyperlinks.Add($ws.Range("J4"), $addrDeJ4, [Type]::Missing, [Type]::Missing, "b7b1e851-a32e-47b1-9c3f-841bf00595e0.md")

Write-Output "Report regenerated for handoff."
